$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (A34): fill previously-empty G5:K5 with numeric values
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -170

# Row 9 (A32): fill previously-empty G9:K9 with numeric values
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = -42.5

# New row 11 (A12) - copy formatting of column A from the last existing data row
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "A12"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 51.089
$ws.Range("H11").Value = 28
$ws.Range("I11").Value = 1100
$ws.Range("J11").Value = 2554.45
$ws.Range("K11").Value = 1100

# New row 12 (A35) - copy formatting of column A from the last existing data row
$ws.Range("A10").Copy($ws.Range("A12"))
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "A35"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 18
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 100
